$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the length for the meeting on row 4 (10:50am - ... -> 10:50am - 11:15)
$ws.Range("B4").Value = "10:50am - 11:15"

# Add the meeting overview for row 4 (new cell E4), matching the style of the
# rest of the row (D4) since E4 previously had no content/format of its own.
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = "CRM, Sprint plan, High Level Design, Project requirements docs"

# Add a new row 5 for the next meeting, carrying over row 4's formatting
# (date number format in A, text style in B:E) before filling in the values.
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)

$ws.Range("A5").Value = "3/29/2018"
$ws.Range("B5").Value = "11:15-11:30"
$ws.Range("C5").Value = "Nicole, Feiyu, Nicole"
$ws.Range("D5").Value = "Jacob"
$ws.Range("E5").Value = "Unit testing, System Testing, Bug list, Performance document "
